$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 130825823
$ws.Cells.Item(2, 2).Value = 57881
$ws.Cells.Item(2, 5).Value = 100049
$ws.Cells.Item(2, 6).Value = "Spillkråka"
$ws.Cells.Item(2, 7).Value = "Dryocopus martius"
$ws.Cells.Item(2, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(2, 13).Value = "äldre spår"
$ws.Cells.Item(2, 16).Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Cells.Item(2, 17).Value = 460947
$ws.Cells.Item(2, 18).Value = 7039711
$ws.Cells.Item(2, 26).Value = "10:38"
$ws.Cells.Item(2, 28).Value = "10:38"
$ws.Cells.Item(2, 29).Value = "Födosökshål på äldre döende gran."

# Row 3
$ws.Cells.Item(3, 1).Value = 130826010
$ws.Cells.Item(3, 2).Value = 91808
$ws.Cells.Item(3, 5).Value = 1202
$ws.Cells.Item(3, 6).Value = "Ullticka"
$ws.Cells.Item(3, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(3, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 16).Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Cells.Item(3, 17).Value = 460971
$ws.Cells.Item(3, 18).Value = 7039688
$ws.Cells.Item(3, 26).Value = "10:47"
$ws.Cells.Item(3, 28).Value = "10:47"
$ws.Cells.Item(3, 29).ClearContents()

# Row 12
$ws.Cells.Item(12, 1).Value = 130826137
$ws.Cells.Item(12, 2).Value = 91808
$ws.Cells.Item(12, 5).Value = 1202
$ws.Cells.Item(12, 6).Value = "Ullticka"
$ws.Cells.Item(12, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(12, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 16).Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Cells.Item(12, 17).Value = 461026
$ws.Cells.Item(12, 18).Value = 7039757
$ws.Cells.Item(12, 26).Value = "10:56"
$ws.Cells.Item(12, 28).Value = "10:56"
$ws.Cells.Item(12, 29).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = 130826287
$ws.Cells.Item(13, 16).Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Cells.Item(13, 17).Value = 461096
$ws.Cells.Item(13, 18).Value = 7039690
$ws.Cells.Item(13, 26).Value = "11:04"
$ws.Cells.Item(13, 28).Value = "11:04"
$ws.Cells.Item(13, 29).Value = "Barkfläkta klenare och grövre granar"

# Row 14
$ws.Cells.Item(14, 1).Value = 130826478
$ws.Cells.Item(14, 2).Value = 57884
$ws.Cells.Item(14, 5).Value = 100109
$ws.Cells.Item(14, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(14, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(14, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(14, 13).Value = "färska spår"
$ws.Cells.Item(14, 16).Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Cells.Item(14, 17).Value = 461220
$ws.Cells.Item(14, 18).Value = 7039590
$ws.Cells.Item(14, 26).Value = "11:25"
$ws.Cells.Item(14, 28).Value = "11:25"
$ws.Cells.Item(14, 29).Value = "Födosök barkfläkt"

# Row 16
$ws.Cells.Item(16, 1).Value = 130825822
$ws.Cells.Item(16, 2).Value = 89193
$ws.Cells.Item(16, 5).Value = 510
$ws.Cells.Item(16, 6).Value = "Doftskinn"
$ws.Cells.Item(16, 7).Value = "Cystostereum murrayi"
$ws.Cells.Item(16, 8).Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Cells.Item(16, 17).Value = 460947
$ws.Cells.Item(16, 18).Value = 7039711
$ws.Cells.Item(16, 26).Value = "10:36"
$ws.Cells.Item(16, 28).Value = "10:36"
$ws.Cells.Item(16, 29).Value = "På granlåga"

# Row 17
$ws.Cells.Item(17, 1).Value = 130826291
$ws.Cells.Item(17, 2).Value = 91808
$ws.Cells.Item(17, 5).Value = 1202
$ws.Cells.Item(17, 6).Value = "Ullticka"
$ws.Cells.Item(17, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(17, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(17, 17).Value = 461106
$ws.Cells.Item(17, 18).Value = 7039672
$ws.Cells.Item(17, 26).Value = "11:04"
$ws.Cells.Item(17, 28).Value = "11:04"
$ws.Cells.Item(17, 29).ClearContents()

# Row 22
$ws.Cells.Item(22, 1).Value = 130826355
$ws.Cells.Item(22, 2).Value = 92535
$ws.Cells.Item(22, 4).Value = "VU"
$ws.Cells.Item(22, 5).Value = 67
$ws.Cells.Item(22, 6).Value = "Sprickporing"
$ws.Cells.Item(22, 7).Value = "Diplomitoporus crustulinus"
$ws.Cells.Item(22, 8).Value = "(Bres.) Domański"
$ws.Cells.Item(22, 16).Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Cells.Item(22, 17).Value = 461117
$ws.Cells.Item(22, 18).Value = 7039629
$ws.Cells.Item(22, 19).Value = 10
$ws.Cells.Item(22, 26).Value = "11:10"
$ws.Cells.Item(22, 28).Value = "11:10"
$ws.Cells.Item(22, 29).Value = "På undersidan av lutande död gran."

# Row 23
$ws.Cells.Item(23, 1).Value = 130826438
$ws.Cells.Item(23, 2).Value = 79243
$ws.Cells.Item(23, 4).Value = "NT"
$ws.Cells.Item(23, 5).Value = 6425
$ws.Cells.Item(23, 6).Value = "Garnlav"
$ws.Cells.Item(23, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(23, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(23, 16).Value = "Brännan, Brännan, Jmt"
$ws.Cells.Item(23, 17).Value = 461220
$ws.Cells.Item(23, 18).Value = 7039590
$ws.Cells.Item(23, 19).Value = 25
$ws.Cells.Item(23, 26).Value = "11:16"
$ws.Cells.Item(23, 28).Value = "11:16"
$ws.Cells.Item(23, 29).Value = "Rikligt i området"
